$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(113).Insert()

$ws.Cells.Item(113, 1).Value = 8
$ws.Cells.Item(113, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(113, 3).Value = "Coquimbo"
$ws.Cells.Item(113, 4).Value = 44438
$ws.Cells.Item(113, 5).Value = 4
$ws.Cells.Item(113, 6).Value = 100114013
$ws.Cells.Item(113, 7).Value = "Zanahoria"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 740
$ws.Cells.Item(113, 11).Value = 5000
$ws.Cells.Item(113, 12).Value = 5500
$ws.Cells.Item(113, 13).Value = 5250
$ws.Cells.Item(113, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(113, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(113, 16).Value = 262
$ws.Cells.Item(113, 17).Value = 20
$ws.Cells.Item(113, 18).Value = "Hortaliza"
